# Updates the cryptos price list (Coin/Link/Price/Volume(1h) columns)
# per commit: "Updated cryptos list on Thu Feb 29 18:48:37 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.689.55"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "3.366.27"
$ws.Range("E3").Value = "  +2.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").Value = "'403.60"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6
$ws.Range("D6").Value = "'126.68"
$ws.Range("E6").Value = "  +13.44%  "

# Row 7
$ws.Range("D7").Value = "'0.605"
$ws.Range("E7").Value = "  +7.36%  "

# Row 8
$ws.Range("E8").Value = "  +0.14%  "

# Row 9
$ws.Range("D9").Value = "'0.666"
$ws.Range("E9").Value = "  +8.41%  "

# Row 10
$ws.Range("E10").Value = "  +15.20%  "

# Row 11
$ws.Range("D11").Value = "'41.84"
$ws.Range("E11").Value = "  +9.23%  "

# Row 12
$ws.Range("E12").Value = "  -0.46%  "

# Row 13
$ws.Range("D13").Value = "3.910.51"
$ws.Range("E13").Value = "  +4.14%  "

# Row 14
$ws.Range("D14").Value = "'8.46"
$ws.Range("E14").Value = "  +4.50%  "

# Row 15
$ws.Range("E15").Value = "  +3.20%  "

# Row 16
$ws.Range("D16").Value = "3.358.14"
$ws.Range("E16").Value = "  +2.04%  "

# Row 17
$ws.Range("D17").Value = "'11.41"
$ws.Range("E17").Value = "  +9.59%  "

# Row 18
$ws.Range("D18").Value = "60.736.61"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19
$ws.Range("E19").Value = "  +2.78%  "

# Row 20
$ws.Range("E20").Value = "  +19.56%  "

# Row 21
$ws.Range("D21").Value = "'3.23"
$ws.Range("E21").Value = "  +1.51%  "

# Row 22
$ws.Range("D22").Value = "'82.01"
$ws.Range("E22").Value = "  +12.47%  "

# Row 23
$ws.Range("D23").Value = "'13.03"
$ws.Range("E23").Value = "  +5.29%  "

# Row 24
$ws.Range("D24").Value = "'304.36"
$ws.Range("E24").Value = "  +3.66%  "

# Row 25
$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  +3.20%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'4.75"
$ws.Range("E26").Value = "  +6.12%  "

# Row 27
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'8.50"
$ws.Range("E27").Value = "  +14.53%  "

# Row 28
$ws.Range("D28").Value = "'29.41"
$ws.Range("E28").Value = "  +3.36%  "

# Row 29
$ws.Range("E29").Value = "  +1.08%  "

# Row 30
$ws.Range("D30").Value = "'0.172"
$ws.Range("E30").Value = "  +1.73%  "

# Row 31
$ws.Range("E31").Value = "  +6.93%  "

# Row 32
$ws.Range("D32").Value = "'11.69"
$ws.Range("E32").Value = "  +5.38%  "

# Row 33
$ws.Range("D33").Value = "'42.50"
$ws.Range("E33").Value = "  +7.48%  "

# Row 34
$ws.Range("D34").Value = "'2.58"
$ws.Range("E34").Value = "  +7.94%  "

# Row 35
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("D36").Value = "'0.0483"
$ws.Range("E36").Value = "  +2.08%  "

# Row 37
$ws.Range("D37").Value = "'52.13"
$ws.Range("E37").Value = "  -0.45%  "

# Row 38
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.16%  "

# Row 39
$ws.Range("E39").Value = "  +3.72%  "

# Row 40
$ws.Range("D40").Value = "'2.95"
$ws.Range("E40").Value = "  -3.65%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'2.05"
$ws.Range("E41").Value = "  +9.88%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.125"
$ws.Range("E42").Value = "  +4.97%  "

# Row 43
$ws.Range("D43").Value = "'135.73"
$ws.Range("E43").Value = "  +0.82%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.91"
$ws.Range("E44").Value = "  +5.23%  "

# Row 45
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'16.82"
$ws.Range("E45").Value = "  +4.70%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.282"
$ws.Range("E46").Value = "  +1.24%  "

# Row 47
$ws.Range("E47").Value = "  +1.52%  "

# Row 48
$ws.Range("D48").Value = "'21.78"
$ws.Range("E48").Value = "  +5.20%  "

# Row 49
$ws.Range("D49").Value = "2.132.83"
$ws.Range("E49").Value = "  +1.32%  "

# Row 50
$ws.Range("D50").Value = "3.705.45"
$ws.Range("E50").Value = "  +2.59%  "

# Row 51
$ws.Range("D51").Value = "'2.35"
$ws.Range("E51").Value = "  +1.15%  "
